$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "emailer and gspread update": the speedtest logger appends another day's
# worth of Date/Time/Download/Upload/Ping samples (columns A:E) below the
# existing log. Find the last populated row in column A (skipping the
# trailing formatting-only row at the bottom of the sheet) and append the
# new "Aug-18" rows right after it.
$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count(), 1).End($xlUp).Row()
$startRow = $lastRow + 1

$newRows = @(
    @("Aug-18", "00:48", 9.57, 10.51, 38.06),
    @("Aug-18", "01:00", 9.26, 10.43, 33.98),
    @("Aug-18", "01:30", 9.59, 10.5, 36.28),
    @("Aug-18", "02:00", 9.94, 10.41, 33.45),
    @("Aug-18", "02:30", 10.37, 10.02, 34.8),
    @("Aug-18", "03:00", 10.33, 9.63, 36.44),
    @("Aug-18", "03:30", 9.96, 9.82, 38.03),
    @("Aug-18", "04:00", 10.36, 10.71, 34.61),
    @("Aug-18", "04:30", 10.35, 10.59, 37.2),
    @("Aug-18", "05:00", 10.35, 10.68, 34.74),
    @("Aug-18", "05:30", 10.33, 10.69, 32.94),
    @("Aug-18", "06:00", 9.73, 10.51, 37.97),
    @("Aug-18", "06:30", 10.31, 10.69, 32.26),
    @("Aug-18", "07:00", 10.36, 9.61, 33.17),
    @("Aug-18", "07:30", 10.17, 10.61, 32.17),
    @("Aug-18", "08:00", 10.33, 10.47, 35.91),
    @("Aug-18", "08:30", 7.27, 8.9, 40.96),
    @("Aug-18", "09:00", 9.53, 9.99, 34.75),
    @("Aug-18", "09:30", 9.47, 4.68, 50.03),
    @("Aug-18", "10:00", 10.31, 10.33, 46.17),
    @("Aug-18", "10:04", 10.3, 9.63, 34.71),
    @("Aug-18", "10:11", 6.0, 10.0, 20.0)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $rec = $newRows[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
}
